# -----------------------------------------------------------------------
# "Arreglos en la documentación"
#
# 1) Adds a new paragraph of text right after the paragraph that ends in
#    "...led y vibración." (end of the MQTT section) and before the
#    "ESP32" heading, plus a fresh empty paragraph after that new text.
# 2) Splits the run that ends "...se alcanzan." (inside the Loop bullet)
#    so a <w:lastRenderedPageBreak/> sits right before "presión".
# 3) Removes the <w:lastRenderedPageBreak/> that currently sits in front
#    of "Funciones Actuadores".
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

function Find-ParagraphIndexEndingWith {
    param([string]$pattern)
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -match $pattern) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# Edit 1: new paragraph about the REST/MQTT publish client, inserted in
# the paragraph that is currently empty right after the "...led y
# vibración." paragraph, followed by a brand new empty paragraph.
# ---------------------------------------------------------------------

$mqttIdx = Find-ParagraphIndexEndingWith "led y vibraci"
$targetIdx = $mqttIdx + 1   # the existing empty paragraph right after it

$p1 = $d.Paragraphs.Item($targetIdx)
$r1 = $p1.Range
$r1.Collapse(1) | Out-Null

$xml1 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Para la implementación del cliente desde el que publicaremos, lo que hemos hecho ha sido crear dentro de la API REST un cliente MQTT que se ha suscrito a los dos canales. Cada vez que hagamos un POST de cualquiera de los actuadores haremos un </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>publish</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">, por lo que se guardará en la BBDD y se actualizará en la placa. De esta forma, nos ahorramos la implementación de MQTT en la </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>app</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>, la cual solo generará peticiones REST, y seguiremos manteniendo MQTT en la placa.</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$r1.InsertXML($xml1)

# ---------------------------------------------------------------------
# Edit 2: split the "...sonido o presión cada 5 segundos..." run so a
# page-break marker sits right before "presión".
# ---------------------------------------------------------------------

$loopIdx = Find-ParagraphIndexEndingWith "valores m.ximos que se alcanzan"
$p2 = $d.Paragraphs.Item($loopIdx)
$r2 = $p2.Range

$xml2 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>
<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Loop</w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>:</w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> En caso de no estar conectados al servidor MQTT nos volvemos a conectar, y luego mediante la función loop() de la librería de MQTT vamos comprobando si existen nuevos datos para los sensores, y luego actualizaremos sus valores. Cada 5 </w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">segundos postearemos </w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">los valores de los sensores, y un segundo después, durante 5 segundos, iremos tomando los valores de sonido y presión y obtendremos el mayor valor obtenido. Esto es porque postear el sonido o </w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>presión cada 5 segundos provocará que nos perdamos varios datos importantes, ya que varían en cuestión de décimas de segundo, y solo necesitamos fijarnos en los valores máximos que se alcanzan.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$r2.InsertXML($xml2)

# ---------------------------------------------------------------------
# Edit 3: drop the <w:lastRenderedPageBreak/> before "Funciones
# Actuadores".
# ---------------------------------------------------------------------

$actIdx = Find-ParagraphIndexEndingWith "^Funciones Actuadores"
$p3 = $d.Paragraphs.Item($actIdx)
$r3 = $p3.Range

$xml3 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>
<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Funciones Actuadores</w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>:</w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Presentamos la función que nos obtiene los IDs de los actuadores a través de una petición HTTP. También tenemos dos funciones análogas entre ellas para escribir los valores recibidos en los actuadores. Estas funciones debemos llamarlas constantemente para conseguir la intermitencia de los actuadores a la intensidad que recibimos.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$r3.InsertXML($xml3)

Write-Host "Edits applied."
